# Fix typo "Back pays" -> "Bank pays" in the Action column of the card data tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Pot luck card data table: row 17, column D ("Savings bond matures, collect £100")
$ws.Range("D17").Value = "Bank pays £100 to the player"

# Opportunity knocks card data table: row 27, column D ("You have won a lip sync battle...")
$ws.Range("D27").Value = "Bank pays player £100"

# Reflect the cell selection / scroll position left by the editor after the edit.
$ws.Activate()
$ws.Range("D28").Select()
$excel.ActiveWindow.ScrollRow = 23
